$d = $word.ActiveDocument

$pairs = @(
    @("37×55=2035", "59×36=2124"),
    @("82×70=5740", "37×96=3552"),
    @("23×12=276", "38×74=2812"),
    @("82×30=2460", "14×48=672"),
    @("20×72=1440", "30×33=990"),
    @("81×56=4536", "65×20=1300"),
    @("84×76=6384", "65×82=5330"),
    @("38×53=2014", "57×94=5358"),
    @("97×16=1552", "16×55=880"),
    @("66×14=924", "61×70=4270"),
    @("20×75=1500", "57×58=3306"),
    @("77×77=5929", "62×69=4278"),
    @("90×32=2880", "89×82=7298"),
    @("24×25=600", "24×58=1392"),
    @("20×79=1580", "86×55=4730"),
    @("84×44=3696", "73×32=2336"),
    @("67×19=1273", "66×18=1188"),
    @("79×88=6952", "96×63=6048"),
    @("70×91=6370", "17×12=204"),
    @("28×79=2212", "44×32=1408"),
    @("72×40=2880", "81×85=6885"),
    @("84×82=6888", "59×79=4661"),
    @("87×34=2958", "68×97=6596"),
    @("98×93=9114", "35×70=2450"),
    @("23×54=1242", "98×17=1666")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
